$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsCtrl  = $wb.Worksheets.Item("BAEPAbCiPC")

# --- BAEPAbCiPC (control-lever) values --------------------------------
# Commit: "3.3.1.1 --> 3.4.0 ... (all fuels affected by production costs)"
# Every fuel row that wasn't already "NOT USED" (gray / style 3) flips from 0 to 1.
$wsCtrl.Range("B3").Value  = 1   # hard coal
$wsCtrl.Range("B4").Value  = 1   # natural gas
$wsCtrl.Range("B9").Value  = 1   # biomass
$wsCtrl.Range("B10").Value = 1   # petroleum gasoline
$wsCtrl.Range("B11").Value = 1   # petroleum diesel
$wsCtrl.Range("B12").Value = 1   # biofuel gasoline
$wsCtrl.Range("B13").Value = 1   # biofuel diesel
$wsCtrl.Range("B14").Value = 1   # jet fuel or kerosene
$wsCtrl.Range("B17").Value = 1   # lignite
$wsCtrl.Range("B18").Value = 1   # crude oil
$wsCtrl.Range("B19").Value = 1   # heavy fuel oil
$wsCtrl.Range("B20").Value = 1   # LPG propane or butane

# Rows 2, 15 and 22 (electricity / heat / hydrogen) previously carried an
# explicit (no-op) "apply fill" style; restore them to the Normal style so
# the now-unused cellXfs entry drops out of the styles table.
$wsCtrl.Range("A2:B2").Style   = "Normal"
$wsCtrl.Range("A15:B15").Style = "Normal"
$wsCtrl.Range("A22:B22").Style = "Normal"

# Column A width nudged slightly narrower.
$wsCtrl.Columns.Item(1).ColumnWidth = 33

# --- View / selection state --------------------------------------------
# Active sheet moves from "About" to "BAEPAbCiPC", and the selection within
# that sheet moves to C16.
$wsCtrl.Activate()
$wsCtrl.Range("C16").Select()
